# Auto-generated edit script: update Kraken_Profits leve profit calculations
# across ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1432.6666
$ws.Range("J17").Value = 1559.2
$ws.Range("L17").Value = 4677.6
$ws.Range("N17").Value = -5013.6
$ws.Range("H125").Value = 806.6667
$ws.Range("I125").Value = 810
$ws.Range("J125").Value = 800
$ws.Range("K125").Value = 7290
$ws.Range("L125").Value = 7200
$ws.Range("M125").Value = -4830
$ws.Range("N125").Value = -12120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -726
$ws.Range("M2").ClearContents()
$ws.Range("H102").Value = 1461.8182
$ws.Range("I102").Value = 1406.9
$ws.Range("J102").Value = 2011
$ws.Range("K102").Value = 1406.9
$ws.Range("L102").Value = 2011
$ws.Range("M102").Value = 215.0999999999999
$ws.Range("N102").Value = -5255
$ws.Range("H116").Value = 500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 500
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 500
$ws.Range("N116").Value = -5088
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 500
$ws.Range("N3").Value = -728
$ws.Range("M3").ClearContents()
$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3346
$ws.Range("H86").Value = 50126.5
$ws.Range("I86").Value = 25249.5
$ws.Range("K86").Value = 25249.5
$ws.Range("M86").Value = -24126.5
$ws.Range("H89").Value = 50126.5
$ws.Range("I89").Value = 25249.5
$ws.Range("K89").Value = 126247.5
$ws.Range("M89").Value = -120631.5
$ws.Range("H99").Value = 866.6667
$ws.Range("I99").Value = 866.6667
$ws.Range("K99").Value = 866.6667
$ws.Range("M99").Value = 631.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 104.666664
$ws.Range("I7").Value = 84.75
$ws.Range("J7").Value = 120.6
$ws.Range("K7").Value = 84.75
$ws.Range("L7").Value = 120.6
$ws.Range("M7").Value = 28.25
$ws.Range("N7").Value = -346.6
$ws.Range("H31").Value = 4300.8237
$ws.Range("I31").Value = 1606.5
$ws.Range("K31").Value = 1606.5
$ws.Range("M31").Value = -1311.5
$ws.Range("H34").Value = 4300.8237
$ws.Range("I34").Value = 1606.5
$ws.Range("K34").Value = 1606.5
$ws.Range("M34").Value = -1404.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1254238
$ws.Range("I7").Value = 4476
$ws.Range("J7").Value = 2504000
$ws.Range("K7").Value = 4476
$ws.Range("L7").Value = 2504000
$ws.Range("M7").Value = -4364
$ws.Range("N7").Value = -2504224
$ws.Range("H8").Value = 1254238
$ws.Range("I8").Value = 4476
$ws.Range("J8").Value = 2504000
$ws.Range("K8").Value = 4476
$ws.Range("L8").Value = 2504000
$ws.Range("M8").Value = -4337
$ws.Range("N8").Value = -2504278
$ws.Range("H9").Value = 1233.3334
$ws.Range("I9").Value = 1100
$ws.Range("K9").Value = 1100
$ws.Range("M9").Value = -930
$ws.Range("H10").Value = 1974.5
$ws.Range("I10").Value = 1975
$ws.Range("J10").Value = 1974
$ws.Range("K10").Value = 1975
$ws.Range("L10").Value = 1974
$ws.Range("M10").Value = -1806
$ws.Range("N10").Value = -2312
$ws.Range("H11").Value = 9401752
$ws.Range("I11").Value = 10335280
$ws.Range("K11").Value = 10335280
$ws.Range("M11").Value = -10335141
$ws.Range("H13").Value = 527.5
$ws.Range("I13").Value = 55
$ws.Range("K13").Value = 55
$ws.Range("M13").Value = 84
$ws.Range("H14").Value = 5077924.5
$ws.Range("J14").Value = 2143074.2
$ws.Range("L14").Value = 2143074.2
$ws.Range("N14").Value = -2143410.2
$ws.Range("H17").Value = 400
$ws.Range("J17").Value = 400
$ws.Range("L17").Value = 400
$ws.Range("N17").Value = -736

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4229.4443
$ws.Range("I40").Value = 4229.4443
$ws.Range("K40").Value = 4229.4443
$ws.Range("M40").Value = -4093.4443
$ws.Range("H61").Value = 3259.6365
$ws.Range("I61").Value = 3320.889
$ws.Range("K61").Value = 3320.889
$ws.Range("M61").Value = -3118.889
$ws.Range("H82").Value = 2005.8
$ws.Range("I82").Value = 1921.6666
$ws.Range("J82").Value = 2132
$ws.Range("K82").Value = 1921.6666
$ws.Range("L82").Value = 2132
$ws.Range("M82").Value = -1560.6666
$ws.Range("N82").Value = -2854
$ws.Range("H85").Value = 2005.8
$ws.Range("I85").Value = 1921.6666
$ws.Range("J85").Value = 2132
$ws.Range("K85").Value = 1921.6666
$ws.Range("L85").Value = 2132
$ws.Range("M85").Value = -673.6666
$ws.Range("N85").Value = -4628
$ws.Range("H100").Value = 4728.8335
$ws.Range("I100").Value = 5386.6
$ws.Range("K100").Value = 5386.6
$ws.Range("M100").Value = -4845.6
$ws.Range("H113").Value = 3259.6365
$ws.Range("I113").Value = 3320.889
$ws.Range("K113").Value = 3320.889
$ws.Range("M113").Value = -1150.889
$ws.Range("H132").Value = 7834.6665
$ws.Range("I132").Value = 7834.6665
$ws.Range("K132").Value = 23503.9995
$ws.Range("M132").Value = -20973.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 3351672
$ws.Range("I52").Value = 4018006.5
$ws.Range("J52").Value = 20000
$ws.Range("K52").Value = 4018006.5
$ws.Range("L52").Value = 20000
$ws.Range("M52").Value = -4017780.5
$ws.Range("N52").Value = -20452
$ws.Range("H100").Value = 8701.888999999999
$ws.Range("I100").Value = 12727.833
$ws.Range("J100").Value = 650
$ws.Range("K100").Value = 25455.666
$ws.Range("L100").Value = 1300
$ws.Range("M100").Value = -24914.666
$ws.Range("N100").Value = -2382
$ws.Range("H107").Value = 2108.8235
$ws.Range("I107").Value = 1590
$ws.Range("K107").Value = 4770
$ws.Range("M107").Value = -2850
$ws.Range("H113").Value = 789.1111
$ws.Range("I113").Value = 914.36365
$ws.Range("J113").Value = 592.2857
$ws.Range("K113").Value = 2743.09095
$ws.Range("L113").Value = 1776.8571
$ws.Range("M113").Value = -573.0909499999998
$ws.Range("N113").Value = -6116.8571
$ws.Range("H122").Value = 2595.7144
$ws.Range("I122").Value = 3234.2
$ws.Range("J122").Value = 999.5
$ws.Range("K122").Value = 9702.599999999999
$ws.Range("L122").Value = 2998.5
$ws.Range("M122").Value = -7252.599999999999
$ws.Range("N122").Value = -7898.5
